# Update the cryptos list (Coin / Link / Price / Volume(1h)) to the latest
# scraped snapshot. Existing rows are updated in place; a few rows were
# re-ordered by the scraper (coins keep their row's rank number in column A,
# but the Coin/Link/Price/Volume values shift to the new rank order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.160.98"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3
$ws.Range("D3").Value = "2.255.77"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'309.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.95%  "

# Row 6
$ws.Range("D6").Value = "'98.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.29%  "

# Row 7
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("E9").Value = "  -3.10%  "

# Row 10
$ws.Range("D10").Value = "'35.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.25%  "

# Row 11
$ws.Range("E11").Value = "  -0.75%  "

# Row 12
$ws.Range("D12").Value = "'7.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.47%  "

# Row 13
$ws.Range("E13").Value = "  -1.65%  "

# Row 14
$ws.Range("D14").Value = "2.598.71"
$ws.Range("E14").Value = "  +0.63%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.841"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.38%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.253.33"
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("D17").Value = "'13.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.85%  "

# Row 18
$ws.Range("D18").Value = "44.037.36"
$ws.Range("E18").Value = "  +1.44%  "

# Row 19
$ws.Range("E19").Value = "  -5.89%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -0.36%  "

# Row 21
$ws.Range("D21").Value = "'6.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.21%  "

# Row 22
$ws.Range("D22").Value = "'65.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "

# Row 23
$ws.Range("D23").Value = "'240.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.05%  "

# Row 24
$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.72%  "

# Row 26
$ws.Range("E26").Value = "  +0.31%  "

# Row 27
$ws.Range("D27").Value = "'10.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "

# Row 28
$ws.Range("D28").Value = "'2.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.68%  "

# Row 29
$ws.Range("D29").Value = "'36.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "

# Row 30
$ws.Range("D30").Value = "'6.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.46%  "

# Row 31
$ws.Range("E31").Value = "  +0.67%  "

# Row 32
$ws.Range("D32").Value = "'157.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.82%  "

# Row 33
$ws.Range("D33").Value = "'3.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.71%  "

# Row 34
$ws.Range("D34").Value = "'0.0829"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.45%  "

# Row 35
$ws.Range("E35").Value = "  +0.21%  "

# Row 36
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.120"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.109"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.67%  "

# Row 38
$ws.Range("D38").Value = "'1.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "

# Row 39
$ws.Range("D39").Value = "'16.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.21%  "

# Row 40
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'3.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.60%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'3.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.46%  "

# Row 42
$ws.Range("E42").Value = "  -3.00%  "

# Row 43
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").Value = "1.768.76"
$ws.Range("E44").Value = "  -1.39%  "

# Row 45
$ws.Range("D45").Value = "'87.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.12%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.194"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.03%  "

# Row 47
$ws.Range("E47").Value = "  -1.32%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'101.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.61%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.26%  "

# Row 50
$ws.Range("D50").Value = "'70.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.88%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'55.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.51%  "

